$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.440.96"
$ws.Range("E2").Value = "  -5.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.970.22"
$ws.Range("E3").Value = "  -7.99%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.05"
$ws.Range("E5").Value = "  -5.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.07"
$ws.Range("E6").Value = "  -10.89%  "

$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.957.22"
$ws.Range("E8").Value = "  -8.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.470"
$ws.Range("E9").Value = "  -13.52%  "

$ws.Range("E10").Value = "  -13.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.96"
$ws.Range("E11").Value = "  -12.66%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.447"
$ws.Range("E12").Value = "  -11.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.71"
$ws.Range("E13").Value = "  -11.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000214"
$ws.Range("E14").Value = "  -12.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.423.43"
$ws.Range("E15").Value = "  -8.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.239.87"
$ws.Range("E16").Value = "  -5.52%  "

$ws.Range("E17").Value = "  -4.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.960.51"
$ws.Range("E18").Value = "  -8.17%  "

$ws.Range("E19").Value = "  -10.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.94"
$ws.Range("E20").Value = "  -14.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.10"
$ws.Range("E21").Value = "  -13.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.648"
$ws.Range("E22").Value = "  -14.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.82"
$ws.Range("E23").Value = "  -11.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.65"
$ws.Range("E24").Value = "  -11.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.03"
$ws.Range("E25").Value = "  -10.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.62"
$ws.Range("E27").Value = "  -16.82%  "

$ws.Range("E28").Value = "  -5.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.31"
$ws.Range("E30").Value = "  -9.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.20"
$ws.Range("E31").Value = "  -14.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.52"
$ws.Range("E32").Value = "  -2.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.07"
$ws.Range("E33").Value = "  -5.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "486.73"
$ws.Range("E34").Value = "  -11.57%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.72"
$ws.Range("E35").Value = "  -3.04%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.09"
$ws.Range("E36").Value = "  -10.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.65"
$ws.Range("E37").Value = "  -13.71%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0387"
$ws.Range("E38").Value = "  -10.62%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0755"
$ws.Range("E39").Value = "  -10.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.114"
$ws.Range("E40").Value = "  -8.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.01"
$ws.Range("E41").Value = "  -12.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.795.61"
$ws.Range("E42").Value = "  -4.64%  "

$ws.Range("E43").Value = "  -0.26%  "

$ws.Range("E44").Value = "  -8.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.229"
$ws.Range("E45").Value = "  -12.09%  "

$ws.Range("B46").Value = "PEPE"
$ws.Range("C46").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0510"
$ws.Range("E46").Value = "  -12.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "115.98"
$ws.Range("E47").Value = "  -5.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.72"
$ws.Range("E48").Value = "  -8.67%  "

$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.93"
$ws.Range("E49").Value = "  -8.07%  "

$ws.Range("E50").Value = "  -8.54%  "

$ws.Range("E51").Value = "  -18.60%  "
